$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (whose D,M,N,O,P,S values move into target row)
$map = @{
  2  = 13
  3  = 5
  4  = 3
  5  = 20
  6  = 10
  7  = 12
  8  = 2
  9  = 18
  10 = 7
  11 = 17
  12 = 4
  13 = 19
  14 = 6
  15 = 15
  16 = 11
  17 = 16
  18 = 9
  19 = 14
  20 = 8
}

# Snapshot original values for the columns that move (D, M, N, O, P, S) for rows 2-20
$cols = @("D", "M", "N", "O", "P", "S")
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
  $rowVals = @{}
  foreach ($c in $cols) {
    $rowVals[$c] = $ws.Range("$c$r").Value2
  }
  $snapshot[$r] = $rowVals
}

# Apply the permutation using the snapshot so rows don't clobber each other
foreach ($targetRow in $map.Keys) {
  $sourceRow = $map[$targetRow]
  $srcVals = $snapshot[$sourceRow]
  foreach ($c in $cols) {
    $ws.Range("$c$targetRow").Value = $srcVals[$c]
  }
}
